$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.042000000000000003
$ws.Range("D3").Value = -1.3149999999999999
$ws.Range("E3").Value = -1.82
$ws.Range("F3").Value = -4.2110000000000003
$ws.Range("G3").Value = -3.5720000000000001

$ws.Range("G4").Select()
